# Update quarterly figures for "Free Cash Flow Margin" (row 19) and
# "Operating Cash Flow Margin" (row 28) on the PINS sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PINS")

# Row 19 - Free Cash Flow Margin
$ws.Range("D19").Value = -0.0394
$ws.Range("E19").Value = 0.0237
$ws.Range("F19").Value = 0.0008
$ws.Range("G19").Value = -0.029

# Row 28 - Operating Cash Flow Margin
$ws.Range("D28").Value = -0.0197
$ws.Range("E28").Value = 0.0508
$ws.Range("F28").Value = 0.0314
$ws.Range("G28").Value = 0.0006
